# .netcore final demo modifications
# Append two new daily-log blocks (45712 / 45713) to the bottom of the
# existing report table on Sheet1, following the same row layout used
# by every earlier day block (Domm / Meeting+Reconsile / General Discussion /
# Study+topic rows / Total row, separated by one blank spacer row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Block 1: 2025-02-24 (serial 45712) -> rows 248-253, row 254 blank ----
$ws.Cells.Item(248,1).Value = 45712
$ws.Cells.Item(248,2).Value = "Domm"
$ws.Cells.Item(248,4).Value = 0.25

$ws.Cells.Item(249,2).Value = "Meeting"
$ws.Cells.Item(249,3).Value = "Reconsile"
$ws.Cells.Item(249,4).Value = 1

$ws.Cells.Item(250,3).Value = "General Discussion"
$ws.Cells.Item(250,4).Value = 0.25

$ws.Cells.Item(251,2).Value = "Study"
$ws.Cells.Item(251,3).Value = "Reconsile Revision & changes"
$ws.Cells.Item(251,4).Value = 1.5

$ws.Cells.Item(252,3).Value = "Editing"
$ws.Cells.Item(252,4).Value = 1.5

$ws.Cells.Item(253,2).Value = "Total"
$ws.Cells.Item(253,4).Formula = "=SUM(D247:D252)"

# ---- Block 2: 2025-02-25 (serial 45713) -> rows 255-261 ----
$ws.Cells.Item(255,1).Value = 45713
$ws.Cells.Item(255,2).Value = "Domm"
$ws.Cells.Item(255,4).Value = 0.25

$ws.Cells.Item(256,2).Value = "Meeting"
$ws.Cells.Item(256,3).Value = "Reconsile"
$ws.Cells.Item(256,4).Value = 0

$ws.Cells.Item(257,3).Value = "General Discussion"
$ws.Cells.Item(257,4).Value = 0.25

$ws.Cells.Item(258,2).Value = "Study"
$ws.Cells.Item(258,3).Value = "Editing"
$ws.Cells.Item(258,4).Value = 3

$ws.Cells.Item(259,3).Value = "Grouping"
$ws.Cells.Item(259,4).Value = 2

$ws.Cells.Item(260,3).Value = "Filtering"
$ws.Cells.Item(260,4).Value = 2.5

$ws.Cells.Item(261,2).Value = "Total"
$ws.Cells.Item(261,4).Formula = "=SUM(D254:D260)"

# ---- Re-apply the recurring block's formatting (borders / number format /
# centering) on top of the values so the new rows look like every other
# day block in the sheet, without disturbing the freshly written values
# or the formula results above. ----
$ws.Range("A240:D245").Copy()
$ws.Range("A248:D253").PasteSpecial(-4122)

$ws.Range("A240:D246").Copy()
$ws.Range("A255:D261").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---- Selection / scroll matches the saved view from the edited workbook ----
[void]$ws.Range("H260").Select()
